$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a "last changed" date serial number.
# Every data row's date is bumped by one day: 45202 (2023-10-03) -> 45203 (2023-10-04).
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) {
    $lastRow = 359
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
